$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$startRow = 59
$startIndicator = 206
$count = 31

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $indicator = $startIndicator + $i
    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($row, 3).Value = "INDICATOR_$indicator"
    $ws.Cells.Item($row, 5).Value = "String"
    $ws.Cells.Item($row, 6).Value = "String"
}

$endRow = $startRow + $count - 1

$ws.Range("A58").Copy()
$ws.Range("A$($startRow):A$($endRow)").PasteSpecial(-4122)
$ws.Range("B58:C58").Copy()
$ws.Range("B$($startRow):C$($endRow)").PasteSpecial(-4122)
$ws.Range("E58:F58").Copy()
$ws.Range("E$($startRow):F$($endRow)").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Activate()
$ws.Range("E59:F89").Select()
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
